# Update attendance summary counts on Sheet1.
# Columns: D=Total Attendance Count, E=Real, F=Duplicate, G=Invalid, H=Absent

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    3  = @{ G = 1; H = 1 }
    4  = @{ D = 1; E = 1 }
    5  = @{ H = 1 }
    6  = @{ D = 1; E = 1 }
    7  = @{ H = 1 }
    8  = @{ H = 1 }
    9  = @{ H = 1 }
    10 = @{ D = 2; E = 1; F = 1 }
    11 = @{ H = 1 }
    12 = @{ D = 1; E = 1 }
    13 = @{ H = 1 }
    14 = @{ D = 1; E = 1 }
    15 = @{ D = 1; E = 1 }
    16 = @{ H = 1 }
    17 = @{ H = 1 }
    18 = @{ D = 1; E = 1 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $cols[$col]
    }
}
